$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 559.8570999999999
$ws.Range("I19").Value = 446.47058
$ws.Range("J19").Value = 735.0909
$ws.Range("K19").Value = 446.47058
$ws.Range("L19").Value = 735.0909
$ws.Range("M19").Value = -271.47058
$ws.Range("N19").Value = -1085.0909

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2492.25
$ws.Range("J51").Value = 2599.6667
$ws.Range("L51").Value = 2599.6667
$ws.Range("N51").Value = -3567.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 351.43332
$ws.Range("I103").Value = 217.73334
$ws.Range("J103").Value = 485.13333
$ws.Range("K103").Value = 653.20002
$ws.Range("L103").Value = 1455.39999
$ws.Range("M103").Value = -67.20001999999999
$ws.Range("N103").Value = -2627.39999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 7576935.5
$ws.Range("J112").Value = 8022625.5
$ws.Range("L112").Value = 24067876.5
$ws.Range("N112").Value = -24070092.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 8624435
$ws.Range("I125").Value = 399.33334
$ws.Range("J125").Value = 11211646
$ws.Range("K125").Value = 3594.00006
$ws.Range("L125").Value = 100904814
$ws.Range("M125").Value = -1134.00006
$ws.Range("N125").Value = -100909734

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1012.73334
$ws.Range("I129").Value = 475.14285
$ws.Range("J129").Value = 1111.7632
$ws.Range("K129").Value = 1425.42855
$ws.Range("L129").Value = 3335.2896
$ws.Range("M129").Value = 3574.57145
$ws.Range("N129").Value = -13335.2896

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 24391496
$ws.Range("I137").Value = 31250852
$ws.Range("K137").Value = 93752556
$ws.Range("M137").Value = -93750006

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 376702.3
$ws.Range("I138").Value = 2060598.1
$ws.Range("J138").Value = 2503.2468
$ws.Range("K138").Value = 6181794.300000001
$ws.Range("L138").Value = 7509.7404
$ws.Range("M138").Value = -6176654.300000001
$ws.Range("N138").Value = -17789.7404

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 26980
$ws.Range("I63").Value = 32975
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 32975
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -32289
$ws.Range("N63").Value = -4372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 26980
$ws.Range("I66").Value = 32975
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 164875
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -161443
$ws.Range("N66").Value = -21864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 22693.5
$ws.Range("J112").Value = 22693.5
$ws.Range("L112").Value = 22693.5
$ws.Range("N112").Value = -25647.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1678.8823
$ws.Range("I122").Value = 1478.5834
$ws.Range("J122").Value = 2159.6
$ws.Range("K122").Value = 4435.7502
$ws.Range("L122").Value = 6478.799999999999
$ws.Range("M122").Value = -1985.7502
$ws.Range("N122").Value = -11378.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2693.5305
$ws.Range("I132").Value = 2300.7437
$ws.Range("J132").Value = 4225.4
$ws.Range("K132").Value = 6902.2311
$ws.Range("L132").Value = 12676.2
$ws.Range("M132").Value = -4372.2311
$ws.Range("N132").Value = -17736.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3781.7273
$ws.Range("I105").Value = 3699.875
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3699.875
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -1952.875
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1123.5333
$ws.Range("I94").Value = 1021.8
$ws.Range("J94").Value = 1174.4
$ws.Range("K94").Value = 1021.8
$ws.Range("L94").Value = 1174.4
$ws.Range("M94").Value = -570.8
$ws.Range("N94").Value = -2076.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12501623
$ws.Range("I99").Value = 31251250
$ws.Range("J99").Value = 1871.3334
$ws.Range("K99").Value = 31251250
$ws.Range("L99").Value = 1871.3334
$ws.Range("M99").Value = -31249752
$ws.Range("N99").Value = -4867.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2625.6
$ws.Range("I122").Value = 1190.5
$ws.Range("J122").Value = 3582.3333
$ws.Range("K122").Value = 3571.5
$ws.Range("L122").Value = 10746.9999
$ws.Range("M122").Value = -1121.5
$ws.Range("N122").Value = -15646.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 12501623
$ws.Range("I126").Value = 31251250
$ws.Range("J126").Value = 1871.3334
$ws.Range("K126").Value = 93753750
$ws.Range("L126").Value = 5614.0002
$ws.Range("M126").Value = -93751280
$ws.Range("N126").Value = -10554.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2917.9092
$ws.Range("I132").Value = 2460.8125
$ws.Range("J132").Value = 4136.8335
$ws.Range("K132").Value = 7382.4375
$ws.Range("L132").Value = 12410.5005
$ws.Range("M132").Value = -4852.4375
$ws.Range("N132").Value = -17470.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2967.6897
$ws.Range("I134").Value = 1331.1666
$ws.Range("J134").Value = 5645.636
$ws.Range("K134").Value = 3993.4998
$ws.Range("L134").Value = 16936.908
$ws.Range("M134").Value = -1458.4998
$ws.Range("N134").Value = -22006.908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1877.7778
$ws.Range("J63").Value = 2000
$ws.Range("L63").Value = 6000
$ws.Range("N63").Value = -7498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2863.7778
$ws.Range("J64").Value = 3259.1428
$ws.Range("L64").Value = 9777.428400000001
$ws.Range("N64").Value = -10317.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 1877.7778
$ws.Range("J66").Value = 2000
$ws.Range("L66").Value = 18000
$ws.Range("N66").Value = -25488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2863.7778
$ws.Range("J67").Value = 3259.1428
$ws.Range("L67").Value = 9777.428400000001
$ws.Range("N67").Value = -11649.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 12656.223
$ws.Range("I87").Value = 4817.6665
$ws.Range("K87").Value = 14452.9995
$ws.Range("M87").Value = -13204.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 5833.3335
$ws.Range("J88").Value = 5833.3335
$ws.Range("L88").Value = 17500.0005
$ws.Range("N88").Value = -18356.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 12656.223
$ws.Range("I90").Value = 4817.6665
$ws.Range("K90").Value = 43358.9985
$ws.Range("M90").Value = -37118.9985

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 5833.3335
$ws.Range("J91").Value = 5833.3335
$ws.Range("L91").Value = 17500.0005
$ws.Range("N91").Value = -20464.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 549.5
$ws.Range("I92").Value = 299
$ws.Range("K92").Value = 897
$ws.Range("M92").Value = 351

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 2671.5789
$ws.Range("J94").Value = 2770
$ws.Range("L94").Value = 8310
$ws.Range("N94").Value = -9662

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 10000
$ws.Range("J104").Value = 10000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -35242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 667.15
$ws.Range("I121").Value = 191.25
$ws.Range("J121").Value = 984.4167
$ws.Range("K121").Value = 573.75
$ws.Range("L121").Value = 2953.2501
$ws.Range("M121").Value = 736.25
$ws.Range("N121").Value = -5573.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1518.9517
$ws.Range("I131").Value = 451.25
$ws.Range("J131").Value = 1677.1296
$ws.Range("K131").Value = 1353.75
$ws.Range("L131").Value = 5031.3888
$ws.Range("M131").Value = 3686.25
$ws.Range("N131").Value = -15111.3888

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2050.7646
$ws.Range("I102").Value = 1774.909
$ws.Range("J102").Value = 2556.5
$ws.Range("K102").Value = 1774.909
$ws.Range("L102").Value = 2556.5
$ws.Range("M102").Value = -152.9090000000001
$ws.Range("N102").Value = -5800.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 32819.5
$ws.Range("J110").Value = 32819.5
$ws.Range("L110").Value = 32819.5
$ws.Range("N110").Value = -40999.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3232.5483
$ws.Range("I136").Value = 2241.0454
$ws.Range("J136").Value = 5656.222
$ws.Range("K136").Value = 6723.1362
$ws.Range("L136").Value = 16968.666
$ws.Range("M136").Value = -4173.1362
$ws.Range("N136").Value = -22068.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3086877
$ws.Range("I107").Value = 3086877
$ws.Range("K107").Value = 9260631
$ws.Range("M107").Value = -9258711

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 50062.43
$ws.Range("I122").Value = 92909.63
$ws.Range("J122").Value = 2930.5
$ws.Range("K122").Value = 278728.89
$ws.Range("L122").Value = 8791.5
$ws.Range("M122").Value = -276278.89
$ws.Range("N122").Value = -13691.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 23813506
$ws.Range("I132").Value = 31254832
$ws.Range("J132").Value = 1260
$ws.Range("K132").Value = 93764496
$ws.Range("L132").Value = 3780
$ws.Range("M132").Value = -93761966
$ws.Range("N132").Value = -8840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7114513.5
$ws.Range("I136").Value = 8798932
$ws.Range("J136").Value = 2521.6667
$ws.Range("K136").Value = 26396796
$ws.Range("L136").Value = 7565.000100000001
$ws.Range("M136").Value = -26394246
$ws.Range("N136").Value = -12665.0001
